$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2843101.8
$ws.Range("I40").Value = 6945833.5
$ws.Range("J40").Value = 2749.2307
$ws.Range("K40").Value = 6945833.5
$ws.Range("L40").Value = 2749.2307
$ws.Range("M40").Value = -6945658.5
$ws.Range("N40").Value = -3099.2307
$ws.Range("H98").Value = 208334990
$ws.Range("I98").Value = 312500480
$ws.Range("J98").Value = 4003
$ws.Range("K98").Value = 312500480
$ws.Range("L98").Value = 4003
$ws.Range("M98").Value = -312498982
$ws.Range("N98").Value = -6999
$ws.Range("H122").Value = 208334990
$ws.Range("I122").Value = 312500480
$ws.Range("J122").Value = 4003
$ws.Range("K122").Value = 937501440
$ws.Range("L122").Value = 12009
$ws.Range("M122").Value = -937498990
$ws.Range("N122").Value = -16909
$ws.Range("H132").Value = 26670654
$ws.Range("I132").Value = 22391486
$ws.Range("J132").Value = 62508676
$ws.Range("K132").Value = 67174458
$ws.Range("L132").Value = 187526028
$ws.Range("M132").Value = -67171928
$ws.Range("N132").Value = -187531088
$ws.Range("H135").Value = 1555.6305
$ws.Range("I135").Value = 1652.9535
$ws.Range("J135").Value = 160.66667
$ws.Range("K135").Value = 14876.5815
$ws.Range("L135").Value = 1446.00003
$ws.Range("M135").Value = -12341.5815
$ws.Range("N135").Value = -6516.00003
$ws.Range("H137").Value = 1375.5
$ws.Range("I137").Value = 1152.24
$ws.Range("K137").Value = 3456.72
$ws.Range("M137").Value = -906.7200000000003
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13902439
$ws.Range("I32").Value = 14936756
$ws.Range("J32").Value = 42599.4
$ws.Range("K32").Value = 14936756
$ws.Range("L32").Value = 42599.4
$ws.Range("M32").Value = -14936469
$ws.Range("N32").Value = -43173.4
$ws.Range("H61").Value = 1677.3
$ws.Range("I61").Value = 1530.6061
$ws.Range("J61").Value = 1856.5927
$ws.Range("K61").Value = 1530.6061
$ws.Range("L61").Value = 1856.5927
$ws.Range("M61").Value = -1318.6061
$ws.Range("N61").Value = -2280.5927
$ws.Range("H122").Value = 3297.3333
$ws.Range("I122").Value = 3655.0833
$ws.Range("J122").Value = 1866.3334
$ws.Range("K122").Value = 10965.2499
$ws.Range("L122").Value = 5599.0002
$ws.Range("M122").Value = -8515.249899999999
$ws.Range("N122").Value = -10499.0002
$ws.Range("H136").Value = 1677.3
$ws.Range("I136").Value = 1530.6061
$ws.Range("J136").Value = 1856.5927
$ws.Range("K136").Value = 4591.8183
$ws.Range("L136").Value = 5569.7781
$ws.Range("M136").Value = -2041.8183
$ws.Range("N136").Value = -10669.7781
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2932004.8
$ws.Range("I134").Value = 7728.5264
$ws.Range("J134").Value = 5856281
$ws.Range("K134").Value = 23185.5792
$ws.Range("L134").Value = 17568843
$ws.Range("M134").Value = -20650.5792
$ws.Range("N134").Value = -17573913
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3532.4722
$ws.Range("I31").Value = 4650.375
$ws.Range("J31").Value = 3213.0715
$ws.Range("K31").Value = 4650.375
$ws.Range("L31").Value = 3213.0715
$ws.Range("M31").Value = -4355.375
$ws.Range("N31").Value = -3803.0715
$ws.Range("H34").Value = 3532.4722
$ws.Range("I34").Value = 4650.375
$ws.Range("J34").Value = 3213.0715
$ws.Range("K34").Value = 4650.375
$ws.Range("L34").Value = 3213.0715
$ws.Range("M34").Value = -4448.375
$ws.Range("N34").Value = -3617.0715
$ws.Range("H58").Value = 41667700
$ws.Range("I58").Value = 66667588
$ws.Range("J58").Value = 1220
$ws.Range("K58").Value = 66667588
$ws.Range("L58").Value = 1220
$ws.Range("M58").Value = -66667385
$ws.Range("N58").Value = -1626
$ws.Range("H59").Value = 31290
$ws.Range("J59").Value = 31290
$ws.Range("L59").Value = 31290
$ws.Range("N59").Value = -33580
$ws.Range("H122").Value = 57693268
$ws.Range("I122").Value = 62500924
$ws.Range("J122").Value = 1400
$ws.Range("K122").Value = 187502772
$ws.Range("L122").Value = 4200
$ws.Range("M122").Value = -187500322
$ws.Range("N122").Value = -9100
$ws.Range("H134").Value = 1498.4884
$ws.Range("I134").Value = 1738.0646
$ws.Range("J134").Value = 879.5833
$ws.Range("K134").Value = 5214.1938
$ws.Range("L134").Value = 2638.7499
$ws.Range("M134").Value = -2679.1938
$ws.Range("N134").Value = -7708.7499
$ws.Range("H136").Value = 41667700
$ws.Range("I136").Value = 66667588
$ws.Range("J136").Value = 1220
$ws.Range("K136").Value = 200002764
$ws.Range("L136").Value = 3660
$ws.Range("M136").Value = -200000214
$ws.Range("N136").Value = -8760
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 17816.5
$ws.Range("I92").Value = 899.6667
$ws.Range("J92").Value = 34733.332
$ws.Range("K92").Value = 2699.0001
$ws.Range("L92").Value = 104199.996
$ws.Range("M92").Value = -1451.0001
$ws.Range("N92").Value = -106695.996
$ws.Range("H126").Value = 1981.4445
$ws.Range("I126").Value = 1260
$ws.Range("J126").Value = 2883.25
$ws.Range("K126").Value = 3780
$ws.Range("L126").Value = 8649.75
$ws.Range("M126").Value = 1160
$ws.Range("N126").Value = -18529.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 20000000
$ws.Range("I40").Value = 20000000
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 20000000
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -19999849
$ws.Range("N40").ClearContents()
$ws.Range("H80").Value = 4547904.5
$ws.Range("I80").Value = 2520.2
$ws.Range("J80").Value = 14288014
$ws.Range("K80").Value = 2520.2
$ws.Range("L80").Value = 14288014
$ws.Range("M80").Value = -1522.2
$ws.Range("N80").Value = -14290010
$ws.Range("H83").Value = 4547904.5
$ws.Range("I83").Value = 2520.2
$ws.Range("J83").Value = 14288014
$ws.Range("K83").Value = 12601
$ws.Range("L83").Value = 71440070
$ws.Range("M83").Value = -7609
$ws.Range("N83").Value = -71450054
$ws.Range("H132").Value = 5469.6
$ws.Range("I132").Value = 1515.579
$ws.Range("J132").Value = 12299.272
$ws.Range("K132").Value = 4546.737
$ws.Range("L132").Value = 36897.81600000001
$ws.Range("M132").Value = -2016.737
$ws.Range("N132").Value = -41957.81600000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2409.818
$ws.Range("I7").Value = 2237.875
$ws.Range("J7").Value = 2868.3333
$ws.Range("K7").Value = 2237.875
$ws.Range("L7").Value = 2868.3333
$ws.Range("M7").Value = -2125.875
$ws.Range("N7").Value = -3092.3333
$ws.Range("H46").Value = 4630218
$ws.Range("I46").Value = 6944977
$ws.Range("K46").Value = 6944977
$ws.Range("M46").Value = -6944789
$ws.Range("H122").Value = 16980.7
$ws.Range("I122").Value = 22092.23
$ws.Range("J122").Value = 7487.857
$ws.Range("K122").Value = 66276.69
$ws.Range("L122").Value = 22463.571
$ws.Range("M122").Value = -63826.69
$ws.Range("N122").Value = -27363.571
$ws.Range("H126").Value = 2409.818
$ws.Range("I126").Value = 2237.875
$ws.Range("J126").Value = 2868.3333
$ws.Range("K126").Value = 6713.625
$ws.Range("L126").Value = 8604.999899999999
$ws.Range("M126").Value = -4243.625
$ws.Range("N126").Value = -13544.9999
$ws.Range("H132").Value = 17246810
$ws.Range("I132").Value = 40002060
$ws.Range("J132").Value = 7984.9395
$ws.Range("K132").Value = 120006180
$ws.Range("L132").Value = 23954.8185
$ws.Range("M132").Value = -120003650
$ws.Range("N132").Value = -29014.8185
$ws.Range("H136").Value = 4399.396
$ws.Range("I136").Value = 3077.0625
$ws.Range("J136").Value = 7044.0625
$ws.Range("K136").Value = 9231.1875
$ws.Range("L136").Value = 21132.1875
$ws.Range("M136").Value = -6681.1875
$ws.Range("N136").Value = -26232.1875
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H115").Value = 40431.11
$ws.Range("J115").Value = 40431.11
$ws.Range("L115").Value = 40431.11
$ws.Range("N115").Value = -43565.11
$ws.Range("H126").Value = 1562.7273
$ws.Range("I126").Value = 998.75
$ws.Range("J126").Value = 3066.6667
$ws.Range("K126").Value = 2996.25
$ws.Range("L126").Value = 9200.000100000001
$ws.Range("M126").Value = -526.25
$ws.Range("N126").Value = -14140.0001
$ws.Range("H132").Value = 12515671
$ws.Range("I132").Value = 18201898
$ws.Range("J132").Value = 5972.88
$ws.Range("K132").Value = 54605694
$ws.Range("L132").Value = 17918.64
$ws.Range("M132").Value = -54603164
$ws.Range("N132").Value = -22978.64
$ws.Range("H136").Value = 3525.647
$ws.Range("I136").Value = 4750.2856
$ws.Range("J136").Value = 2034.7826
$ws.Range("K136").Value = 14250.8568
$ws.Range("L136").Value = 6104.3478
$ws.Range("M136").Value = -11700.8568
$ws.Range("N136").Value = -11204.3478

Write-Host "Applied 226 cell updates across 8 sheets"